$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Numeric "Qty executed upto date" column (C) updates ---
$ws.Range("C8").Value = 99
$ws.Range("C9").Value = 14
$ws.Range("C10").Value = 17
$ws.Range("C11").Value = 26
$ws.Range("C12").Value = 36
$ws.Range("C13").Value = 53
$ws.Range("C14").Value = 32
$ws.Range("C15").Value = 59
$ws.Range("C16").Value = 26
$ws.Range("C17").Value = 12

# --- "Upto date Amount" column (G) text values, recomputed as Qty * Rate ---
# These cells are stored as text (e.g. "3584.00"), so force text format
# before assigning so Excel keeps them as strings instead of numbers.
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "3584.00"

$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "8024.00"

$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "17212.00"

$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "7208.00"

$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "736.00"

# --- Grand total rows (19 and 21), columns G and H ---
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "36764.00"

$ws.Range("H19").NumberFormat = "@"
$ws.Range("H19").Value = "36764.00"

$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "36764.00"

$ws.Range("H21").NumberFormat = "@"
$ws.Range("H21").Value = "36764.00"
